$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11: GuilID / object field, mirroring the existing rows' layout.
$ws.Range("A11").Value = "GuilID"
$ws.Range("B11").Value = "object"
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = $true
$ws.Range("E11").Value = $true
$ws.Range("F11").Value = $true
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = "Friend"
$ws.Range("J11").Value = "工会ID"

# Match the style used by the other data rows (s="1" -> text-formatted cells).
$ws.Range("A11:B11").NumberFormat = "@"
$ws.Range("I11:J11").NumberFormat = "@"

$ws.Range("E19").Select()
